$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AddOpportunity")
$ws2 = $wb.Worksheets.Item("Users")

# --- AddOpportunity sheet: add the two new "CSDN" rows first so the new
#     shared-string entries land in the same order the source workbook has
#     them (CSDN-0000001546, Eric Winthrop, DC, Dental, HC - Healthcare). ---
$ws1.Range("E6").Value = "CSDN-0000001546"
$ws1.Range("E7").Value = "CSDN-0000001546"

# Staff column (AF) - new contact replaces both Ayati Arvind / Aaron Schultz
$ws1.Range("AF2").Value = "Eric Winthrop"
$ws1.Range("AF3").Value = "Eric Winthrop"

# PrimaryOffice column (K) - AM -> DC
$ws1.Range("K2").Value = "DC"
$ws1.Range("K3").Value = "DC"

# Sector column (E) - Dealership & Rental Services -> Dental
$ws1.Range("E2").Value = "Dental"
$ws1.Range("E3").Value = "Dental"

# IndustryGroup column (D) - BUS - Business Services -> HC - Healthcare
$ws1.Range("D2").Value = "HC - Healthcare"
$ws1.Range("D3").Value = "HC - Healthcare"

# New formatting (wrap text + vertical center) on the Sector cells and the
# two new CSDN cells - build the style once on E2 then copy/paste it onto
# the other three cells so only a single new cellXf is produced.
$ws1.Range("E2").VerticalAlignment = -4108
$ws1.Range("E2").WrapText = $true
$ws1.Range("E2").Copy()
$ws1.Range("E3").PasteSpecial(-4122)
$ws1.Range("E6").PasteSpecial(-4122)
$ws1.Range("E7").PasteSpecial(-4122)

# --- Users sheet: same contact swap ---
$ws2.Range("A2").Value = "Eric Winthrop"

# --- Selection / active cell bookkeeping, matching the saved view state ---
# Set the non-active sheet's remembered selection first...
$ws2.Range("A2").Select()
# ...then land back on AddOpportunity so it stays the active tab.
$ws1.Range("D2:D3").Select()
